# Apply the "season record" columns (Wins, Losses, Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): AD1=Wins, AE1=Losses, AF1=Ties ---
# Copy the existing header style (from AC1, the last used header cell) onto
# the new header cells so they keep the bold/centered/bordered formatting.
$headerSrc = $ws.Range("AC1")
$headerDst = $ws.Range("AD1:AF1")
$headerSrc.Copy($headerDst)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-52): AD=64 (Wins), AE=97 (Losses), AF=0 (Ties) ---
for ($r = 2; $r -le 52; $r++) {
    $ws.Cells.Item($r, 30).Value = 64   # column AD
    $ws.Cells.Item($r, 31).Value = 97   # column AE
    $ws.Cells.Item($r, 32).Value = 0    # column AF
}
